# Weekly update: insert a new week of data (rows 237-238) above the
# existing rows, shifting the previous rows down by two (old 237-240
# become new 239-242).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 237.
$ws.Rows("237:238").Insert()

# --- Row 237: Betarraga, Primera, new week (2022-02-03) ---
$ws.Range("A237").Value = 1
$ws.Range("B237").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C237").Value = "Arica y Parinacota"
$ws.Range("D237").Value = 44595
$ws.Range("E237").Value = 15
$ws.Range("F237").Value = 100114014
$ws.Range("G237").Value = "Betarraga"
$ws.Range("H237").Value = "Sin especificar"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 1200
$ws.Range("K237").Value = 450
$ws.Range("L237").Value = 500
$ws.Range("M237").Value = 475
$ws.Range("N237").Value = "`$/paquete 4 unidades"
$ws.Range("O237").Value = "Región de Arica y Parinacota"
$ws.Range("P237").Value = 119
$ws.Range("Q237").Value = 4
$ws.Range("R237").Value = "Hortaliza"

# --- Row 238: Betarraga, Segunda, new week (2022-02-03) ---
$ws.Range("A238").Value = 1
$ws.Range("B238").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C238").Value = "Arica y Parinacota"
$ws.Range("D238").Value = 44595
$ws.Range("E238").Value = 15
$ws.Range("F238").Value = 100114014
$ws.Range("G238").Value = "Betarraga"
$ws.Range("H238").Value = "Sin especificar"
$ws.Range("I238").Value = "Segunda"
$ws.Range("J238").Value = 1200
$ws.Range("K238").Value = 450
$ws.Range("L238").Value = 500
$ws.Range("M238").Value = 475
$ws.Range("N238").Value = "`$/paquete 5 unidades"
$ws.Range("O238").Value = "Región de Arica y Parinacota"
$ws.Range("P238").Value = 95
$ws.Range("Q238").Value = 5
$ws.Range("R238").Value = "Hortaliza"
